$d = $word.ActiveDocument

# Start at the end of the last paragraph ("Mit terrain.party ... importiert")
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)  # wdCollapseEnd

# --- New bullet: Terrain hinzugefügt ... ---
$tail.InsertParagraphAfter()
$tail.Collapse(0)
$tail.InsertAfter("Terrain hinzugefügt und durch die Tools selber bearbeitet (height, noise, smooth height)")

# --- New bullet: Assets Outdoor Ground Textures, Conifers, und Grass Flowers Pack Free ---
$tail2 = $d.Paragraphs.Last.Range
$tail2.Collapse(0)
$tail2.InsertParagraphAfter()
$tail2.Collapse(0)
$tail2.InsertAfter("Assets Outdoor Ground Textures,")

$tail3 = $d.Paragraphs.Last.Range
$tail3.Collapse(0)
$tail3.InsertAfter(" Conifers, und Grass Flowers Pack Free")
